$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D1").Value = "joinable_relation.id"
$ws.Range("E1").Value = "joinable_relation.foreign_field"
$ws.Range("F1").Value = "joinable_relation.another_foreign_field"
$ws.Range("G1").Value = "joinable_relation.nested_joinable_relation.foreign_field"
